# Generate Report for Handoff
#
# b.md is now ready to be handed off again (it is no longer "in sync" /
# a content duplicate of a.md): update the Overview sheet plus the two
# per-locale sheets (zh-cn, de-de) with the new status, new handoff
# artifact names/timestamps, and the "stale handback" error message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-27 14:37:39"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"

# "False" would be auto-coerced to a real Boolean by plain assignment;
# force it to stay literal text (matching the original shared-string
# cell type) via the text leading-quote, then drop the resulting
# "quote prefix" style so the cell's style index is unaffected.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"

$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-27 14:37:35"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4dad23bb083439d75dfd36af07ad54dfe5be834/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b2b9910c25173d6ae69f0fba2f659b10b74187d/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet - row 3 is the "b.md" row
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"

$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"

$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-27 14:37:39"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4dad23bb083439d75dfd36af07ad54dfe5be834/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b2b9910c25173d6ae69f0fba2f659b10b74187d/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 40
